# model3_df_results.xlsx - "new README - requires cleaning"
#
# Refreshes the R^2 / RMSE / U metric columns (C:E, rows 2-9) with new
# values, and re-paints the RMSE (D) and U (E) cell backgrounds with the
# refreshed green color-scale (light text on the dark-green end, dark text
# on the light-green end), matching the new numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's OLE_COLOR is 0x00BBGGRR (little-endian RGB), not the more common
# 0xRRGGBB - build it from a hex "RRGGBB" string ourselves since the usual
# RGB() helper isn't available in this host.
# (NOTE: avoid using $r anywhere else in this script - it collides with the
# helper's local $r and silently clobbers it in this host's flat scoping.)
function Hex-Color([string]$hexColor) {
    $rVal = [Convert]::ToInt32($hexColor.Substring(0, 2), 16)
    $gVal = [Convert]::ToInt32($hexColor.Substring(2, 2), 16)
    $bVal = [Convert]::ToInt32($hexColor.Substring(4, 2), 16)
    return $rVal + ($gVal * 256) + ($bVal * 65536)
}

$rowData = @(
    @{ Row = 2; C = 0.3518; D = 1.0065; DFill = "006428"; DFont = "F1F1F1"; E = 2.1292; EFill = "0E7936"; EFont = "F1F1F1" },
    @{ Row = 3; C = 0.3758; D = 0.9819; DFill = "00441B"; DFont = "F1F1F1"; E = 1.9681; EFill = "00441B"; EFont = "F1F1F1" },
    @{ Row = 4; C = 0.3439; D = 1.0052; DFill = "006328"; DFont = "F1F1F1"; E = 2.0183; EFill = "005622"; EFont = "F1F1F1" },
    @{ Row = 5; C = 0.3368; D = 1.0151; DFill = "026F2E"; DFont = "F1F1F1"; E = 2.0538; EFill = "006328"; EFont = "F1F1F1" },
    @{ Row = 6; C = 0.2439; D = 1.0847; DFill = "52B365"; DFont = "000000"; E = 2.2083; EFill = "268E47"; EFont = "000000" },
    @{ Row = 7; C = 0.1728; D = 1.1322; DFill = "9BD696"; DFont = "000000"; E = 2.6692; EFill = "CCEBC6"; EFont = "000000" },
    @{ Row = 8; C = 0.1075; D = 1.1841; DFill = "D7EFD1"; DFont = "000000"; E = 2.7723; EFill = "E7F6E2"; EFont = "000000" },
    @{ Row = 9; C = 0.0514; D = 1.2296; DFill = "F7FCF5"; DFont = "000000"; E = 2.8770; EFill = "F7FCF5"; EFont = "000000" }
)

foreach ($entry in $rowData) {
    $n = $entry.Row

    $ws.Range("C$n").Value = $entry.C

    $dCell = $ws.Range("D$n")
    $dCell.Value = $entry.D
    $dCell.Interior.Color = Hex-Color $entry.DFill
    $dCell.Font.Color = Hex-Color $entry.DFont

    $eCell = $ws.Range("E$n")
    $eCell.Value = $entry.E
    $eCell.Interior.Color = Hex-Color $entry.EFill
    $eCell.Font.Color = Hex-Color $entry.EFont
}
